$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("derek.monahan@hotmail.com", "151Fa04124@4517"),
    @("madelene.crooks@hotmail.com", "151Fa04124@4517"),
    @("josefina.weissnat@hotmail.com", "151Fa04124@4517"),
    @("lester.marvin@yahoo.com", "151Fa04124@4517"),
    @("blair.langosh@yahoo.com", "ula20xmk8"),
    @("giovanni.larkin@yahoo.com", "xo3dhmd3vh"),
    @("wynell.tremblay@gmail.com", "wrww8cmt"),
    @("winston.feil@gmail.com", "yndn2uh4"),
    @("shawn.wehner@hotmail.com", "ldksakt9"),
    @("ferdinand.tromp@gmail.com", "rtm1s1fur"),
    @("micha.mertz@hotmail.com", "4uhbevikl0ndaz"),
    @("blanche.gulgowski@yahoo.com", "gzyycon69tdv8s"),
    @("erlinda.padberg@hotmail.com", "kgijjmw44"),
    @("adolfo.terry@yahoo.com", "los70114gy"),
    @("mohammad.mosciski@yahoo.com", "lyfrzaqjpl"),
    @("ferdinand.kris@yahoo.com", "xvgdmq69t9qbk"),
    @("emory.towne@yahoo.com", "p9yrtmidqhnt"),
    @("dion.bruen@gmail.com", "hqmpbjn84u"),
    @("malisa.hammes@hotmail.com", "nlfnial0"),
    @("patrice.considine@yahoo.com", "l6ednspb6k9w8xy"),
    @("bonnie.crooks@yahoo.com", "b0gvuqkzd"),
    @("perry.brakus@gmail.com", "7iuqcpkz0g6"),
    @("domenica.beahan@yahoo.com", "f2ilwmlnplyf53"),
    @("antione.schowalter@gmail.com", "54vg99my4f5l"),
    @("joesph.ebert@gmail.com", "x3ommy81y4qe26")
)

$row = 9
foreach ($pair in $data) {
    $ws.Cells.Item($row, 1).Value = $pair[0]
    $ws.Cells.Item($row, 2).Value = $pair[1]
    $row = $row + 1
}

$ws.Range("I16").Select()
